$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 73.5
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 33

# Row 3
$ws.Range("C3").Value = 92.7
$ws.Range("D3").Value = 48
$ws.Range("E3").Value = 44

# Row 4
$ws.Range("C4").Value = 65.09999999999999
$ws.Range("D4").Value = 47
$ws.Range("E4").Value = 30

# Row 5
$ws.Range("C5").Value = 71
$ws.Range("D5").Value = 43
$ws.Range("E5").Value = 30

# Row 6
$ws.Range("C6").Value = 88.90000000000001
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 44
